$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.886.74"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").Value = "2.586.73"
$ws.Range("E3").Value = "  +2.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'302.99"
$ws.Range("E5").Value = "  +1.87%  "

# Row 6
$ws.Range("D6").Value = "'96.97"
$ws.Range("E6").Value = "  +3.81%  "

# Row 7
$ws.Range("D7").Value = "'0.575"
$ws.Range("E7").Value = "  +0.85%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +0.65%  "

# Row 10
$ws.Range("D10").Value = "'36.58"
$ws.Range("E10").Value = "  +0.79%  "

# Row 11
$ws.Range("E11").Value = "  +0.92%  "

# Row 12
$ws.Range("D12").Value = "'7.71"
$ws.Range("E12").Value = "  +1.83%  "

# Row 13
$ws.Range("E13").Value = "  +7.00%  "

# Row 14
$ws.Range("D14").Value = "2.567.47"
$ws.Range("E14").Value = "  +1.83%  "

# Row 15
$ws.Range("D15").Value = "'0.884"
$ws.Range("E15").Value = "  +2.04%  "

# Row 16
$ws.Range("D16").Value = "'14.39"
$ws.Range("E16").Value = "  +2.18%  "

# Row 17
$ws.Range("D17").Value = "42.919.59"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18
$ws.Range("D18").Value = "'12.91"
$ws.Range("E18").Value = "  +4.92%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0992"
$ws.Range("E19").Value = "  +2.95%  "

# Row 20
$ws.Range("E20").Value = "  +1.97%  "

# Row 21
$ws.Range("D21").Value = "'72.00"
$ws.Range("E21").Value = "  -0.82%  "

# Row 22
$ws.Range("D22").Value = "'254.79"
$ws.Range("E22").Value = "  -1.58%  "

# Row 23
$ws.Range("E23").Value = "  +2.54%  "

# Row 24
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  -0.98%  "

# Row 25
$ws.Range("D25").Value = "'28.68"
$ws.Range("E25").Value = "  -1.21%  "

# Row 26
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("E27").Value = "  +2.56%  "

# Row 28
$ws.Range("E28").Value = "  +7.13%  "

# Row 29
$ws.Range("E29").Value = "  -3.88%  "

# Row 30
$ws.Range("D30").Value = "'6.05"
$ws.Range("E30").Value = "  +1.43%  "

# Row 31
$ws.Range("D31").Value = "'155.70"
$ws.Range("E31").Value = "  +2.71%  "

# Row 32
$ws.Range("D32").Value = "'2.18"
$ws.Range("E32").Value = "  -0.71%  "

# Row 33
$ws.Range("E33").Value = "  -0.19%  "

# Row 34
$ws.Range("E34").Value = "  +1.70%  "

# Row 35
$ws.Range("D35").Value = "'3.37"
$ws.Range("E35").Value = "  -3.01%  "

# Row 36
$ws.Range("D36").Value = "'18.41"
$ws.Range("E36").Value = "  +11.93%  "

# Row 37
$ws.Range("E37").Value = "  +0.63%  "

# Row 38
$ws.Range("E38").Value = "  +1.18%  "

# Row 39
$ws.Range("D39").Value = "'23.43"
$ws.Range("E39").Value = "  -1.81%  "

# Row 40
$ws.Range("B40").Value = "ApeXProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D40").Value = "'2.08"
$ws.Range("E40").Value = "  +29.87%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'3.90"
$ws.Range("E41").Value = "  +1.74%  "

# Row 42
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "'3.42"
$ws.Range("E42").Value = "  -0.92%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0311"
$ws.Range("E43").Value = "  +1.03%  "

# Row 44
$ws.Range("D44").Value = "2.070.32"
$ws.Range("E44").Value = "  +2.59%  "

# Row 45
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.31%  "

# Row 46
$ws.Range("D46").Value = "'9.24"
$ws.Range("E46").Value = "  +4.21%  "

# Row 47
$ws.Range("D47").Value = "'85.24"
$ws.Range("E47").Value = "  -0.54%  "

# Row 48
$ws.Range("D48").Value = "'76.63"
$ws.Range("E48").Value = "  +12.05%  "

# Row 49
$ws.Range("D49").Value = "2.835.26"
$ws.Range("E49").Value = "  +2.60%  "

# Row 50
$ws.Range("D50").Value = "'106.36"
$ws.Range("E50").Value = "  +3.42%  "

# Row 51
$ws.Range("D51").Value = "'1.68"
$ws.Range("E51").Value = "  +2.42%  "
